$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused columns D through N, leaving only A:C of data.
$ws.Range("D1:N6").Clear()

# New header rows (row 1 = field names, row 2 = field types).
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "icon"
$ws.Range("C1").Value = "res"

$ws.Range("A2").Value = "key"
$ws.Range("B2").Value = "string"
$ws.Range("C2").Value = "string"

# Data rows: simple sequential numbers 1-4 in every column.
for ($r = 3; $r -le 6; $r++) {
    $n = $r - 2
    $ws.Cells.Item($r, 1).Value = $n
    $ws.Cells.Item($r, 2).Value = $n
    $ws.Cells.Item($r, 3).Value = $n
}

# Matches the cursor position left behind in the saved file.
$null = $ws.Range("E12").Select()
